$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cells to Text format first so numeric-looking price strings
# (e.g. "213.37", "0.0621") are preserved exactly as text, matching
# the original inlineStr cell type, instead of being auto-converted
# to Excel numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.683.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.37'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.87%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0621'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.14'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0835'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.864.50'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.631.13'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.93%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.678.35'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.04'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.36'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.55%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.39'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.61%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.93'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.32'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.93%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.42'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.69'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0516'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.84%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.96'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.93%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.67%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.165.61'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0167'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.808'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.20%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.502'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.36'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.774.06'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.44'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +10.64%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.47%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.86%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.56'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.30%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.410'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.72%  '
